# First set of edits after R&R
# - Rename row label "Forced commitment" (A5) -> "Mandatory structured"
# - Rename row label "Choice commitment" (A7) -> "Choice " (trailing space)
# - Move the active selection to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("repeat_loans")

$ws.Range("A5").Value = "Mandatory structured"
$ws.Range("A7").Value = "Choice "

$ws.Range("A7").Select()
